$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = 2.3
$ws.Range("A2").Value = 1.1
$ws.Range("A3").Value = 0.1
$ws.Range("A4").Value = 3.3
$ws.Range("A5").Value = 2.05
$ws.Range("A6").Value = 1.975
$ws.Range("A7").Value = 0.6000000000000001
$ws.Range("A8").Value = -0
$ws.Range("A9").Value = 3.1
$ws.Range("A10").Value = 2.8
